# feat: add 2022-Q3 data
#
# 1) "总计" summary sheet: insert a new row right under the header with the
#    2022-Q3 totals, pushing the existing quarters down by one row.
# 2) Insert a brand-new "2022-Q3" worksheet (cloned from "2022-Q2" so it
#    keeps the same header/row styling) right before the existing
#    "2022-Q2" sheet, then overwrite its single data row with the Q3 figures.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update the "总计" (totals) sheet
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")

# Make room for the new quarter just below the header row.
$totals.Rows.Item(2).Insert()

# Row 2 (freshly inserted) inherited row-1's header formatting; restore the
# plain data-row look by copying A3's format (already correct, since row 3
# is the former row 2 that got pushed down) onto the new A2, and clearing
# the formatting on B2:D2 back to the workbook default.
$totals.Cells.Item(3, 1).Copy()
$totals.Cells.Item(2, 1).PasteSpecial(-4122)
$totals.Range("B2:D2").ClearFormats()

$totals.Cells.Item(2, 1).Value = 0
$totals.Cells.Item(2, 2).Value = "2022-Q3"
$totals.Cells.Item(2, 3).Value = 1
$totals.Cells.Item(2, 4).Value = 1.14

# ---------------------------------------------------------------------
# 2) Insert the new "2022-Q3" worksheet
# ---------------------------------------------------------------------
$oldQ2 = $wb.Worksheets.Item("2022-Q2")

# Cloning keeps the header/border/bold styling identical to the other
# quarter sheets without having to rebuild it cell by cell. Excel names the
# freshly-inserted clone "2022-Q2 (2)" (the original keeps "2022-Q2"), so
# grab it by that auto-generated name before renaming it to "2022-Q3".
$oldQ2.Copy($oldQ2)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# B2/C2 (fund code/name) stay the same; refresh the Q3 figures. The values
# are quoted with a leading apostrophe so they stay text (matching the
# source data's string-typed numeric columns) instead of being coerced to
# numbers.
$q3.Cells.Item(2, 4).Value = "'35.14"
$q3.Cells.Item(2, 5).Value = "'85.39"
$q3.Cells.Item(2, 6).Value = "'3.24"
$q3.Cells.Item(2, 7).Value = "'1.1385"
$q3.Cells.Item(2, 8).Value = 9
